$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 903
$ws.Range("F4").Value = 4423
$ws.Range("F6").Value = 428
$ws.Range("F7").Value = 3528
$ws.Range("F8").Value = 1002
$ws.Range("F11").Value = 320
$ws.Range("F12").Value = 316
$ws.Range("F13").Value = 2398
$ws.Range("F14").Value = 1260
$ws.Range("F17").Value = 9
$ws.Range("F18").Value = 524
$ws.Range("F19").Value = 251
$ws.Range("F21").Value = 9952
$ws.Range("F22").Value = 5947
$ws.Range("F24").Value = 204
$ws.Range("F25").Value = 819
$ws.Range("F26").Value = 137
$ws.Range("F27").Value = 836
$ws.Range("F28").Value = 3525
$ws.Range("F30").Value = 964
$ws.Range("F31").Value = 457
$ws.Range("F32").Value = 111
$ws.Range("F33").Value = 238
$ws.Range("F35").Value = 221
$ws.Range("F36").Value = 4802
$ws.Range("F38").Value = 1073
$ws.Range("F39").Value = 144
$ws.Range("F40").Value = 18
$ws.Range("F41").Value = 55
$ws.Range("F42").Value = 473

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F12").Value = 127
$ws.Range("F15").Value = 3522
$ws.Range("F16").Value = 76

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 8689
$ws.Range("F3").Value = 416
$ws.Range("F4").Value = 1553

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 8690
$ws.Range("F3").Value = 903
$ws.Range("F4").Value = 1553
$ws.Range("F6").Value = 4423
$ws.Range("F8").Value = 428
$ws.Range("F9").Value = 3528
$ws.Range("F10").Value = 1002
$ws.Range("F13").Value = 2399
$ws.Range("F18").Value = 1260
$ws.Range("F21").Value = 127
$ws.Range("F22").Value = 524
$ws.Range("F23").Value = 251
$ws.Range("F25").Value = 9952
$ws.Range("F26").Value = 3522
$ws.Range("F27").Value = 76
$ws.Range("F29").Value = 204
$ws.Range("F30").Value = 819
$ws.Range("F31").Value = 137
$ws.Range("F32").Value = 836
$ws.Range("F33").Value = 3525
$ws.Range("F35").Value = 964
$ws.Range("F36").Value = 457
$ws.Range("F37").Value = 111
$ws.Range("F40").Value = 221
$ws.Range("F41").Value = 4802
$ws.Range("F42").Value = 1073
$ws.Range("F44").Value = 55
$ws.Range("F45").Value = 473
